# "Update countries & provincias Spain"
#
# The source COVID dashboard data was refreshed: Spain's (and several
# other countries') case counts changed, which shuffled the ranking of a
# handful of neighbouring rows (Kuwait/Hong Kong/Armenia/Azerbaiyan,
# Libano/Principado de Andorra, Etiopia/Jamaica/Congo, and
# Malaui/Montserrat/Republica de Africa Central/Islas Turcas y Caicos),
# and the "updated at" timestamp moved from 11:22 to 11:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 11:52"

# España (row 5)
$ws.Range("B5").Value = 157022
$ws.Range("C5").Value = 3800
$ws.Range("D5").Value = 55668
$ws.Range("E5").Value = 85511
$ws.Range("G5").Value = 396
$ws.Range("H5").Value = 15843

# Japon (row 14)
$ws.Range("B14").Value = 24172
$ws.Range("C14").Value = 121
$ws.Range("E14").Value = 12614
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 958

# Panama (row 22)
$ws.Range("B22").Value = 10095
$ws.Range("C22").Value = 127
$ws.Range("D22").Value = 1061
$ws.Range("E22").Value = 8942
$ws.Range("F22").Value = 164
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 92

# Croacia (row 37)
$ws.Range("B37").Value = 4346
$ws.Range("C37").Value = 118
$ws.Range("D37").Value = 1830
$ws.Range("E37").Value = 2446
$ws.Range("F37").Value = 69
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 70

# Rows 69-72 shuffle: Kuwait moves up ahead of Hong Kong / Armenia / Azerbaiyan
$ws.Range("A69").Value = "Kuwait"
$ws.Range("B69").Value = 993
$ws.Range("C69").Value = 83
$ws.Range("D69").Value = 123
$ws.Range("E69").Value = 869
$ws.Range("F69").Value = 26
$ws.Range("H69").Value = 1

$ws.Range("A70").Value = "Hong Kong"
$ws.Range("B70").Value = 974
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 293
$ws.Range("E70").Value = 677
$ws.Range("F70").Value = 14
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 4

$ws.Range("A71").Value = "Armenia"
$ws.Range("B71").Value = 937
$ws.Range("C71").Value = 16
$ws.Range("D71").Value = 149
$ws.Range("E71").Value = 777
$ws.Range("F71").Value = 30
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 11

$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 926
$ws.Range("D72").Value = 101
$ws.Range("E72").Value = 816
$ws.Range("F72").Value = 27
$ws.Range("H72").Value = 9

# Rows 84-85 swap: Libano moves ahead of Principado de Andorra
$ws.Range("A84").Value = "Libano"
$ws.Range("B84").Value = 609
$ws.Range("C84").Value = 27
$ws.Range("D84").Value = 67
$ws.Range("E84").Value = 522
$ws.Range("F84").Value = 28
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 20

$ws.Range("A85").Value = "Principado de Andorra"
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 58
$ws.Range("E85").Value = 500
$ws.Range("F85").Value = 17
$ws.Range("H85").Value = 25

# Butan (row 124)
$ws.Range("B124").Value = 136
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 99
$ws.Range("E124").Value = 36

# Rows 140-142 shuffle: Etiopia moves up ahead of Jamaica / Congo
$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 65
$ws.Range("C140").Value = 9
$ws.Range("D140").Value = 4
$ws.Range("E140").Value = 59
$ws.Range("F140").Value = 2
$ws.Range("H140").Value = 2

$ws.Range("A141").Value = "Jamaica"
$ws.Range("B141").Value = 63
$ws.Range("D141").Value = 13
$ws.Range("E141").Value = 46
$ws.Range("H141").Value = 4

$ws.Range("A142").Value = "Congo"
$ws.Range("B142").Value = 60
$ws.Range("D142").Value = 5
$ws.Range("F142").Value = 0
$ws.Range("H142").Value = 5

# Rows 193-196 shuffle: Malaui moves up ahead of Montserrat / Republica de
# Africa Central / Islas Turcas y Caicos
$ws.Range("A193").Value = "Malaui"
$ws.Range("C193").Value = 1
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = 1
$ws.Range("H193").Value = 1

$ws.Range("A194").Value = "Montserrat"
$ws.Range("B194").Value = 9
$ws.Range("E194").Value = 7
$ws.Range("H194").Value = 2

$ws.Range("A195").Value = "Republica de Africa Central"
$ws.Range("E195").Value = 8
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Islas Turcas y Caicos"
$ws.Range("F196").Value = 0
